# Fasl-Fas.xlsx update: refreshed TPM-derived NATMI values, added Resolving-Mac sending-cluster rows (17-21)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "FAPs"
$ws.Cells.Item(2,2).Value = "Fasl"
$ws.Cells.Item(2,3).Value = "Fas"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 0.1330356666666667
$ws.Cells.Item(2,8).Value = 0.399107
$ws.Cells.Item(2,9).Value = 0.2921233835506897
$ws.Cells.Item(2,10).Value = 0.3075293673451303
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 4.364279333333333
$ws.Cells.Item(2,14).Value = 13.092838
$ws.Cells.Item(2,15).Value = 0.1623236474417427
$ws.Cells.Item(2,16).Value = 0.1693585125158563
$ws.Cells.Item(2,17).Value = 0.5806048106295555
$ws.Cells.Item(2,18).Value = 5.225443295666
$ws.Cells.Item(2,19).Value = 0.04741853312097115
$ws.Cells.Item(2,20).Value = 0.05208271620851362

# Row 3
$ws.Cells.Item(3,1).Value = "FAPs"
$ws.Cells.Item(3,2).Value = "Fasl"
$ws.Cells.Item(3,3).Value = "Fas"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 0.6666666666666666
$ws.Cells.Item(3,7).Value = 0.1330356666666667
$ws.Cells.Item(3,8).Value = 0.399107
$ws.Cells.Item(3,9).Value = 0.2921233835506897
$ws.Cells.Item(3,10).Value = 0.3075293673451303
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 8.811908
$ws.Cells.Item(3,14).Value = 26.435724
$ws.Cells.Item(3,15).Value = 0.327747364050729
$ws.Cells.Item(3,16).Value = 0.341951446578635
$ws.Cells.Item(3,17).Value = 1.172298055385333
$ws.Cells.Item(3,18).Value = 10.550682498468
$ws.Cells.Item(3,19).Value = 0.09574266893631864
$ws.Cells.Item(3,20).Value = 0.1051601120290798

# Row 4
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Fasl"
$ws.Cells.Item(4,3).Value = "Fas"
$ws.Cells.Item(4,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 0.6666666666666666
$ws.Cells.Item(4,7).Value = 0.1330356666666667
$ws.Cells.Item(4,8).Value = 0.399107
$ws.Cells.Item(4,9).Value = 0.2921233835506897
$ws.Cells.Item(4,10).Value = 0.3075293673451303
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 5.303883333333334
$ws.Cells.Item(4,14).Value = 15.91165
$ws.Cells.Item(4,15).Value = 0.1972709862305183
$ws.Cells.Item(4,16).Value = 0.205820416908307
$ws.Cells.Item(4,17).Value = 0.7056056551722223
$ws.Cells.Item(4,18).Value = 6.350450896550001
$ws.Cells.Item(4,19).Value = 0.05762746797404052
$ws.Cells.Item(4,20).Value = 0.06329582259852262

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Fasl"
$ws.Cells.Item(5,3).Value = "Fas"
$ws.Cells.Item(5,4).Value = "MuSCs"
$ws.Cells.Item(5,5).Value = 2
$ws.Cells.Item(5,6).Value = 0.6666666666666666
$ws.Cells.Item(5,7).Value = 0.1330356666666667
$ws.Cells.Item(5,8).Value = 0.399107
$ws.Cells.Item(5,9).Value = 0.2921233835506897
$ws.Cells.Item(5,10).Value = 0.3075293673451303
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 3.3504315
$ws.Cells.Item(5,14).Value = 6.700863
$ws.Cells.Item(5,15).Value = 0.1246149066192622
$ws.Cells.Item(5,16).Value = 0.08667702069272822
$ws.Cells.Item(5,17).Value = 0.4457268882235
$ws.Cells.Item(5,18).Value = 2.674361329341
$ws.Cells.Item(5,19).Value = 0.03640292816247212
$ws.Cells.Item(5,20).Value = 0.02665572933699548

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Fasl"
$ws.Cells.Item(6,3).Value = "Fas"
$ws.Cells.Item(6,4).Value = "Resolving-Mac"
$ws.Cells.Item(6,5).Value = 2
$ws.Cells.Item(6,6).Value = 0.6666666666666666
$ws.Cells.Item(6,7).Value = 0.1330356666666667
$ws.Cells.Item(6,8).Value = 0.399107
$ws.Cells.Item(6,9).Value = 0.2921233835506897
$ws.Cells.Item(6,10).Value = 0.3075293673451303
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 5.055779666666667
$ws.Cells.Item(6,14).Value = 15.167339
$ws.Cells.Item(6,15).Value = 0.1880430956577478
$ws.Cells.Item(6,16).Value = 0.1961926033044734
$ws.Cells.Item(6,17).Value = 0.6725990184747778
$ws.Cells.Item(6,18).Value = 6.053391166273
$ws.Cells.Item(6,19).Value = 0.05493178535688729
$ws.Cells.Item(6,20).Value = 0.06033498717201884

# Row 7
$ws.Cells.Item(7,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(7,2).Value = "Fasl"
$ws.Cells.Item(7,3).Value = "Fas"
$ws.Cells.Item(7,4).Value = "ECs"
$ws.Cells.Item(7,5).Value = 1
$ws.Cells.Item(7,6).Value = 0.3333333333333333
$ws.Cells.Item(7,7).Value = 0.1356786666666667
$ws.Cells.Item(7,8).Value = 0.407036
$ws.Cells.Item(7,9).Value = 0.2979269557961613
$ws.Cells.Item(7,10).Value = 0.3136390080020958
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 4.364279333333333
$ws.Cells.Item(7,14).Value = 13.092838
$ws.Cells.Item(7,15).Value = 0.1623236474417427
$ws.Cells.Item(7,16).Value = 0.1693585125158563
$ws.Cells.Item(7,17).Value = 0.5921396009075556
$ws.Cells.Item(7,18).Value = 5.329256408168001
$ws.Cells.Item(7,19).Value = 0.04836059013604776
$ws.Cells.Item(7,20).Value = 0.0531174358621837

# Row 8
$ws.Cells.Item(8,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(8,2).Value = "Fasl"
$ws.Cells.Item(8,3).Value = "Fas"
$ws.Cells.Item(8,4).Value = "FAPs"
$ws.Cells.Item(8,5).Value = 1
$ws.Cells.Item(8,6).Value = 0.3333333333333333
$ws.Cells.Item(8,7).Value = 0.1356786666666667
$ws.Cells.Item(8,8).Value = 0.407036
$ws.Cells.Item(8,9).Value = 0.2979269557961613
$ws.Cells.Item(8,10).Value = 0.3136390080020958
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 8.811908
$ws.Cells.Item(8,14).Value = 26.435724
$ws.Cells.Item(8,15).Value = 0.327747364050729
$ws.Cells.Item(8,16).Value = 0.341951446578635
$ws.Cells.Item(8,17).Value = 1.195587928229333
$ws.Cells.Item(8,18).Value = 10.760291354064
$ws.Cells.Item(8,19).Value = 0.09764477444184991
$ws.Cells.Item(8,20).Value = 0.1072493124898048

# Row 9
$ws.Cells.Item(9,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(9,2).Value = "Fasl"
$ws.Cells.Item(9,3).Value = "Fas"
$ws.Cells.Item(9,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(9,5).Value = 1
$ws.Cells.Item(9,6).Value = 0.3333333333333333
$ws.Cells.Item(9,7).Value = 0.1356786666666667
$ws.Cells.Item(9,8).Value = 0.407036
$ws.Cells.Item(9,9).Value = 0.2979269557961613
$ws.Cells.Item(9,10).Value = 0.3136390080020958
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 5.303883333333334
$ws.Cells.Item(9,14).Value = 15.91165
$ws.Cells.Item(9,15).Value = 0.1972709862305183
$ws.Cells.Item(9,16).Value = 0.205820416908307
$ws.Cells.Item(9,17).Value = 0.7196238188222224
$ws.Cells.Item(9,18).Value = 6.476614369400001
$ws.Cells.Item(9,19).Value = 0.05877234439456477
$ws.Cells.Item(9,20).Value = 0.06455331138569921

# Row 10
$ws.Cells.Item(10,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(10,2).Value = "Fasl"
$ws.Cells.Item(10,3).Value = "Fas"
$ws.Cells.Item(10,4).Value = "MuSCs"
$ws.Cells.Item(10,5).Value = 1
$ws.Cells.Item(10,6).Value = 0.3333333333333333
$ws.Cells.Item(10,7).Value = 0.1356786666666667
$ws.Cells.Item(10,8).Value = 0.407036
$ws.Cells.Item(10,9).Value = 0.2979269557961613
$ws.Cells.Item(10,10).Value = 0.3136390080020958
$ws.Cells.Item(10,11).Value = 2
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 3.3504315
$ws.Cells.Item(10,14).Value = 6.700863
$ws.Cells.Item(10,15).Value = 0.1246149066192622
$ws.Cells.Item(10,16).Value = 0.08667702069272822
$ws.Cells.Item(10,17).Value = 0.454582078678
$ws.Cells.Item(10,18).Value = 2.727492472068
$ws.Cells.Item(10,19).Value = 0.0371261397758997
$ws.Cells.Item(10,20).Value = 0.02718529478664441

# Row 11
$ws.Cells.Item(11,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(11,2).Value = "Fasl"
$ws.Cells.Item(11,3).Value = "Fas"
$ws.Cells.Item(11,4).Value = "Resolving-Mac"
$ws.Cells.Item(11,5).Value = 1
$ws.Cells.Item(11,6).Value = 0.3333333333333333
$ws.Cells.Item(11,7).Value = 0.1356786666666667
$ws.Cells.Item(11,8).Value = 0.407036
$ws.Cells.Item(11,9).Value = 0.2979269557961613
$ws.Cells.Item(11,10).Value = 0.3136390080020958
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 5.055779666666667
$ws.Cells.Item(11,14).Value = 15.167339
$ws.Cells.Item(11,15).Value = 0.1880430956577478
$ws.Cells.Item(11,16).Value = 0.1961926033044734
$ws.Cells.Item(11,17).Value = 0.6859614441337778
$ws.Cells.Item(11,18).Value = 6.173652997204
$ws.Cells.Item(11,19).Value = 0.05602310704779916
$ws.Cells.Item(11,20).Value = 0.06153365347776375

# Row 12
$ws.Cells.Item(12,1).Value = "MuSCs"
$ws.Cells.Item(12,2).Value = "Fasl"
$ws.Cells.Item(12,3).Value = "Fas"
$ws.Cells.Item(12,4).Value = "ECs"
$ws.Cells.Item(12,5).Value = 1
$ws.Cells.Item(12,6).Value = 0.5
$ws.Cells.Item(12,7).Value = 0.0684425
$ws.Cells.Item(12,8).Value = 0.136885
$ws.Cells.Item(12,9).Value = 0.1502879278890229
$ws.Cells.Item(12,10).Value = 0.1054758684990195
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 4.364279333333333
$ws.Cells.Item(12,14).Value = 13.092838
$ws.Cells.Item(12,15).Value = 0.1623236474417427
$ws.Cells.Item(12,16).Value = 0.1693585125158563
$ws.Cells.Item(12,17).Value = 0.2987021882716667
$ws.Cells.Item(12,18).Value = 1.79221312963
$ws.Cells.Item(12,19).Value = 0.0243952846214078
$ws.Cells.Item(12,20).Value = 0.017863236195312

# Row 13
$ws.Cells.Item(13,1).Value = "MuSCs"
$ws.Cells.Item(13,2).Value = "Fasl"
$ws.Cells.Item(13,3).Value = "Fas"
$ws.Cells.Item(13,4).Value = "FAPs"
$ws.Cells.Item(13,5).Value = 1
$ws.Cells.Item(13,6).Value = 0.5
$ws.Cells.Item(13,7).Value = 0.0684425
$ws.Cells.Item(13,8).Value = 0.136885
$ws.Cells.Item(13,9).Value = 0.1502879278890229
$ws.Cells.Item(13,10).Value = 0.1054758684990195
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 8.811908
$ws.Cells.Item(13,14).Value = 26.435724
$ws.Cells.Item(13,15).Value = 0.327747364050729
$ws.Cells.Item(13,16).Value = 0.341951446578635
$ws.Cells.Item(13,17).Value = 0.6031090132900001
$ws.Cells.Item(13,18).Value = 3.61865407974
$ws.Cells.Item(13,19).Value = 0.04925647221427328
$ws.Cells.Item(13,20).Value = 0.03606762581237759

# Row 14
$ws.Cells.Item(14,1).Value = "MuSCs"
$ws.Cells.Item(14,2).Value = "Fasl"
$ws.Cells.Item(14,3).Value = "Fas"
$ws.Cells.Item(14,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(14,5).Value = 1
$ws.Cells.Item(14,6).Value = 0.5
$ws.Cells.Item(14,7).Value = 0.0684425
$ws.Cells.Item(14,8).Value = 0.136885
$ws.Cells.Item(14,9).Value = 0.1502879278890229
$ws.Cells.Item(14,10).Value = 0.1054758684990195
$ws.Cells.Item(14,11).Value = 3
$ws.Cells.Item(14,12).Value = 1
$ws.Cells.Item(14,13).Value = 5.303883333333334
$ws.Cells.Item(14,14).Value = 15.91165
$ws.Cells.Item(14,15).Value = 0.1972709862305183
$ws.Cells.Item(14,16).Value = 0.205820416908307
$ws.Cells.Item(14,17).Value = 0.3630110350416668
$ws.Cells.Item(14,18).Value = 2.17806621025
$ws.Cells.Item(14,19).Value = 0.02964744775320855
$ws.Cells.Item(14,20).Value = 0.02170908722823396

# Row 15
$ws.Cells.Item(15,1).Value = "MuSCs"
$ws.Cells.Item(15,2).Value = "Fasl"
$ws.Cells.Item(15,3).Value = "Fas"
$ws.Cells.Item(15,4).Value = "MuSCs"
$ws.Cells.Item(15,5).Value = 1
$ws.Cells.Item(15,6).Value = 0.5
$ws.Cells.Item(15,7).Value = 0.0684425
$ws.Cells.Item(15,8).Value = 0.136885
$ws.Cells.Item(15,9).Value = 0.1502879278890229
$ws.Cells.Item(15,10).Value = 0.1054758684990195
$ws.Cells.Item(15,11).Value = 2
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 3.3504315
$ws.Cells.Item(15,14).Value = 6.700863
$ws.Cells.Item(15,15).Value = 0.1246149066192622
$ws.Cells.Item(15,16).Value = 0.08667702069272822
$ws.Cells.Item(15,17).Value = 0.22931190793875
$ws.Cells.Item(15,18).Value = 0.917247631755
$ws.Cells.Item(15,19).Value = 0.018728116099893
$ws.Cells.Item(15,20).Value = 0.009142334036472992

# Row 16
$ws.Cells.Item(16,1).Value = "MuSCs"
$ws.Cells.Item(16,2).Value = "Fasl"
$ws.Cells.Item(16,3).Value = "Fas"
$ws.Cells.Item(16,4).Value = "Resolving-Mac"
$ws.Cells.Item(16,5).Value = 1
$ws.Cells.Item(16,6).Value = 0.5
$ws.Cells.Item(16,7).Value = 0.0684425
$ws.Cells.Item(16,8).Value = 0.136885
$ws.Cells.Item(16,9).Value = 0.1502879278890229
$ws.Cells.Item(16,10).Value = 0.1054758684990195
$ws.Cells.Item(16,11).Value = 3
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = 5.055779666666667
$ws.Cells.Item(16,14).Value = 15.167339
$ws.Cells.Item(16,15).Value = 0.1880430956577478
$ws.Cells.Item(16,16).Value = 0.1961926033044734
$ws.Cells.Item(16,17).Value = 0.3460301998358334
$ws.Cells.Item(16,18).Value = 2.076181199015
$ws.Cells.Item(16,19).Value = 0.02826060720024023
$ws.Cells.Item(16,20).Value = 0.02069358522662294

# Row 17
$ws.Cells.Item(17,1).Value = "Resolving-Mac"
$ws.Cells.Item(17,2).Value = "Fasl"
$ws.Cells.Item(17,3).Value = "Fas"
$ws.Cells.Item(17,4).Value = "ECs"
$ws.Cells.Item(17,5).Value = 1
$ws.Cells.Item(17,6).Value = 0.3333333333333333
$ws.Cells.Item(17,7).Value = 0.1182523333333333
$ws.Cells.Item(17,8).Value = 0.354757
$ws.Cells.Item(17,9).Value = 0.259661732764126
$ws.Cells.Item(17,10).Value = 0.2733557561537542
$ws.Cells.Item(17,11).Value = 3
$ws.Cells.Item(17,12).Value = 1
$ws.Cells.Item(17,13).Value = 4.364279333333333
$ws.Cells.Item(17,14).Value = 13.092838
$ws.Cells.Item(17,15).Value = 0.1623236474417427
$ws.Cells.Item(17,16).Value = 0.1693585125158563
$ws.Cells.Item(17,17).Value = 0.516086214485111
$ws.Cells.Item(17,18).Value = 4.644775930366
$ws.Cells.Item(17,19).Value = 0.04214923956331601
$ws.Cells.Item(17,20).Value = 0.04629512424984695

# Row 18
$ws.Cells.Item(18,1).Value = "Resolving-Mac"
$ws.Cells.Item(18,2).Value = "Fasl"
$ws.Cells.Item(18,3).Value = "Fas"
$ws.Cells.Item(18,4).Value = "FAPs"
$ws.Cells.Item(18,5).Value = 1
$ws.Cells.Item(18,6).Value = 0.3333333333333333
$ws.Cells.Item(18,7).Value = 0.1182523333333333
$ws.Cells.Item(18,8).Value = 0.354757
$ws.Cells.Item(18,9).Value = 0.259661732764126
$ws.Cells.Item(18,10).Value = 0.2733557561537542
$ws.Cells.Item(18,11).Value = 3
$ws.Cells.Item(18,12).Value = 1
$ws.Cells.Item(18,13).Value = 8.811908
$ws.Cells.Item(18,14).Value = 26.435724
$ws.Cells.Item(18,15).Value = 0.327747364050729
$ws.Cells.Item(18,16).Value = 0.341951446578635
$ws.Cells.Item(18,17).Value = 1.042028682118667
$ws.Cells.Item(18,18).Value = 9.378258139068
$ws.Cells.Item(18,19).Value = 0.0851034484582871
$ws.Cells.Item(18,20).Value = 0.09347439624737287

# Row 19
$ws.Cells.Item(19,1).Value = "Resolving-Mac"
$ws.Cells.Item(19,2).Value = "Fasl"
$ws.Cells.Item(19,3).Value = "Fas"
$ws.Cells.Item(19,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(19,5).Value = 1
$ws.Cells.Item(19,6).Value = 0.3333333333333333
$ws.Cells.Item(19,7).Value = 0.1182523333333333
$ws.Cells.Item(19,8).Value = 0.354757
$ws.Cells.Item(19,9).Value = 0.259661732764126
$ws.Cells.Item(19,10).Value = 0.2733557561537542
$ws.Cells.Item(19,11).Value = 3
$ws.Cells.Item(19,12).Value = 1
$ws.Cells.Item(19,13).Value = 5.303883333333334
$ws.Cells.Item(19,14).Value = 15.91165
$ws.Cells.Item(19,15).Value = 0.1972709862305183
$ws.Cells.Item(19,16).Value = 0.205820416908307
$ws.Cells.Item(19,17).Value = 0.6271965798944446
$ws.Cells.Item(19,18).Value = 5.644769219050001
$ws.Cells.Item(19,19).Value = 0.05122372610870442
$ws.Cells.Item(19,20).Value = 0.05626219569585122

# Row 20
$ws.Cells.Item(20,1).Value = "Resolving-Mac"
$ws.Cells.Item(20,2).Value = "Fasl"
$ws.Cells.Item(20,3).Value = "Fas"
$ws.Cells.Item(20,4).Value = "MuSCs"
$ws.Cells.Item(20,5).Value = 1
$ws.Cells.Item(20,6).Value = 0.3333333333333333
$ws.Cells.Item(20,7).Value = 0.1182523333333333
$ws.Cells.Item(20,8).Value = 0.354757
$ws.Cells.Item(20,9).Value = 0.259661732764126
$ws.Cells.Item(20,10).Value = 0.2733557561537542
$ws.Cells.Item(20,11).Value = 2
$ws.Cells.Item(20,12).Value = 1
$ws.Cells.Item(20,13).Value = 3.3504315
$ws.Cells.Item(20,14).Value = 6.700863
$ws.Cells.Item(20,15).Value = 0.1246149066192622
$ws.Cells.Item(20,16).Value = 0.08667702069272822
$ws.Cells.Item(20,17).Value = 0.3961963425485
$ws.Cells.Item(20,18).Value = 2.377178055291
$ws.Cells.Item(20,19).Value = 0.03235772258099739
$ws.Cells.Item(20,20).Value = 0.02369366253261532

# Row 21
$ws.Cells.Item(21,1).Value = "Resolving-Mac"
$ws.Cells.Item(21,2).Value = "Fasl"
$ws.Cells.Item(21,3).Value = "Fas"
$ws.Cells.Item(21,4).Value = "Resolving-Mac"
$ws.Cells.Item(21,5).Value = 1
$ws.Cells.Item(21,6).Value = 0.3333333333333333
$ws.Cells.Item(21,7).Value = 0.1182523333333333
$ws.Cells.Item(21,8).Value = 0.354757
$ws.Cells.Item(21,9).Value = 0.259661732764126
$ws.Cells.Item(21,10).Value = 0.2733557561537542
$ws.Cells.Item(21,11).Value = 3
$ws.Cells.Item(21,12).Value = 1
$ws.Cells.Item(21,13).Value = 5.055779666666667
$ws.Cells.Item(21,14).Value = 15.167339
$ws.Cells.Item(21,15).Value = 0.1880430956577478
$ws.Cells.Item(21,16).Value = 0.1961926033044734
$ws.Cells.Item(21,17).Value = 0.5978577424025556
$ws.Cells.Item(21,18).Value = 5.380719681623
$ws.Cells.Item(21,19).Value = 0.0488275960528211
$ws.Cells.Item(21,20).Value = 0.05363037742806788

